$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Log hours worked + task description ("Metodología" entry) for 22-25 Apr (rows 59-62)
$ws.Range("E59").Value = 2
$ws.Range("F59").Value = " ●Investigar/avanzar TFG en overleaf: Metodología"

$ws.Range("E60").Value = 2
$ws.Range("F60").Value = " ●Investigar/avanzar TFG en overleaf: Metodología"

$ws.Range("E61").Value = 1
$ws.Range("F61").Value = " ●Investigar/avanzar TFG en overleaf: Metodología"

$ws.Range("E62").Value = 1
$ws.Range("F62").Value = " ●Investigar/avanzar TFG en overleaf: Metodología"

# Extend the gray "current day" highlight band (as used on rows 55-56) down onto
# rows 63-64, by copying their formatting across.
$ws.Range("C55:F56").Copy() | Out-Null
$ws.Range("C63:F64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection on the last cell touched in this session
$ws.Range("F67").Select() | Out-Null

$wb.Save()
